# Adiciona a caixa "setor" (Rectangle 4) ao slide 2, logo apos a caixa de
# data ("06/05/2025"), duplicando a caixa "resumo" (Rectangle 4) ja
# existente para herdar toda a formatacao (fonte, cor, paragrafos, locks),
# e entao ajustando posicao, texto e altura para os valores finais.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Localiza a caixa "Rectangle 4" (id 44, texto "resumo") que serve de
# modelo para o novo rotulo "setor".
$src = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Id -eq 44) {
        $src = $cand
        break
    }
}

$dup = $src.Duplicate()
$newShape = $dup.Item(1)

# Posiciona a nova caixa (valores convertidos de EMU para pontos: /12700).
$newShape.Left = -5.215196850393701
$newShape.Top = 321.4892125984252

# Primeiro paragrafo recebe o texto "setor"; o "`r" preserva o segundo
# paragrafo vazio (mesma estrutura do modelo).
$newShape.TextFrame.TextRange.Text = "setor`r"

# Restaura a altura original da caixa apos o recalculo do autofit.
$newShape.Height = 83.4672440944882
